$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Apr 08 17:54:25 EDT 2024"
$ws.Range("B3").Value = "Mon Apr 08 17:54:36 EDT 2024"
$ws.Range("B4").Value = "Mon Apr 08 17:54:47 EDT 2024"
$ws.Range("B5").Value = "Mon Apr 08 17:54:59 EDT 2024"
$ws.Range("B6").Value = "Mon Apr 08 17:55:10 EDT 2024"
$ws.Range("B7").Value = "Mon Apr 08 17:55:21 EDT 2024"
$ws.Range("B8").Value = "Mon Apr 08 17:55:33 EDT 2024"
$ws.Range("B9").Value = "Mon Apr 08 17:55:44 EDT 2024"
$ws.Range("B10").Value = "Mon Apr 08 17:55:56 EDT 2024"
$ws.Range("B11").Value = "Mon Apr 08 17:56:07 EDT 2024"
$ws.Range("B12").Value = "Mon Apr 08 17:56:18 EDT 2024"
$ws.Range("B13").Value = "Mon Apr 08 17:56:29 EDT 2024"
$ws.Range("B14").Value = "Mon Apr 08 17:56:41 EDT 2024"
$ws.Range("B15").Value = "Mon Apr 08 17:56:53 EDT 2024"
$ws.Range("B16").Value = "Mon Apr 08 17:57:05 EDT 2024"
$ws.Range("B17").Value = "Mon Apr 08 17:57:17 EDT 2024"
$ws.Range("B18").Value = "Mon Apr 08 17:57:29 EDT 2024"
$ws.Range("B19").Value = "Mon Apr 08 17:57:41 EDT 2024"
$ws.Range("B20").Value = "Mon Apr 08 17:57:53 EDT 2024"
$ws.Range("B21").Value = "Mon Apr 08 17:58:05 EDT 2024"
$ws.Range("B22").Value = "Mon Apr 08 17:58:17 EDT 2024"
$ws.Range("B23").Value = "Mon Apr 08 17:58:29 EDT 2024"
$ws.Range("B24").Value = "Mon Apr 08 17:58:41 EDT 2024"
$ws.Range("B25").Value = "Mon Apr 08 17:58:52 EDT 2024"
$ws.Range("B26").Value = "Mon Apr 08 17:59:04 EDT 2024"
$ws.Range("B27").Value = "Mon Apr 08 17:59:16 EDT 2024"
$ws.Range("B28").Value = "Mon Apr 08 17:59:30 EDT 2024"
$ws.Range("B29").Value = "Mon Apr 08 17:59:42 EDT 2024"
$ws.Range("B30").Value = "Mon Apr 08 17:59:54 EDT 2024"
$ws.Range("B31").Value = "Mon Apr 08 18:00:05 EDT 2024"
$ws.Range("B32").Value = "Mon Apr 08 18:00:17 EDT 2024"
$ws.Range("B33").Value = "Mon Apr 08 18:00:29 EDT 2024"
$ws.Range("B34").Value = "Mon Apr 08 18:00:41 EDT 2024"
$ws.Range("B35").Value = "Mon Apr 08 18:00:53 EDT 2024"
$ws.Range("B36").Value = "Mon Apr 08 18:01:04 EDT 2024"
$ws.Range("B37").Value = "Mon Apr 08 18:01:15 EDT 2024"
$ws.Range("B38").Value = "Mon Apr 08 18:01:27 EDT 2024"
$ws.Range("B39").Value = "Mon Apr 08 18:01:39 EDT 2024"
$ws.Range("B40").Value = "Mon Apr 08 18:01:50 EDT 2024"
$ws.Range("B41").Value = "Mon Apr 08 18:02:02 EDT 2024"
$ws.Range("B42").Value = "Mon Apr 08 18:02:13 EDT 2024"
$ws.Range("B43").Value = "Mon Apr 08 18:02:25 EDT 2024"
$ws.Range("B44").Value = "Mon Apr 08 18:02:37 EDT 2024"
$ws.Range("B45").Value = "Mon Apr 08 18:02:49 EDT 2024"
$ws.Range("B46").Value = "Mon Apr 08 18:03:01 EDT 2024"
$ws.Range("B47").Value = "Mon Apr 08 18:03:14 EDT 2024"
$ws.Range("B48").Value = "Mon Apr 08 18:03:26 EDT 2024"
$ws.Range("B49").Value = "Mon Apr 08 18:03:38 EDT 2024"
$ws.Range("B50").Value = "Mon Apr 08 18:03:50 EDT 2024"
$ws.Range("B51").Value = "Mon Apr 08 18:04:03 EDT 2024"
$ws.Range("B52").Value = "Mon Apr 08 18:04:15 EDT 2024"
$ws.Range("B53").Value = "Mon Apr 08 18:04:27 EDT 2024"
$ws.Range("B54").Value = "Mon Apr 08 18:04:39 EDT 2024"
